$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 (index 10): "@IBActionfunc" -> "@IBAction func "
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(4)
$tr10 = $shp10.TextFrame.TextRange
$run = $tr10.Find("@IBActionfunc")
$run.Text = "@IBAction func "

# ---------------------------------------------------------------------------
# Slide 6 (index 6): "iflet" -> "if let"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(4)
$tr6 = $shp6.TextFrame.TextRange
$run = $tr6.Find("iflet")
$run.Text = "if let"

# ---------------------------------------------------------------------------
# Slide 8 (index 8): split the comment run so the placeholder name
# "John Appleseed" becomes "Nem Sothea", keeping the same green run
# formatting around it (three runs instead of one).
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(4)
$tr8 = $shp8.TextFrame.TextRange
$found8 = $tr8.Find('// "John Appleseed" (Force-unwrap - careful!)')
$start8 = $found8.Start
$found8.Text = '// "Nem Sothea" (Force-unwrap - careful!)'

# Re-select just the new name ("// `"" is 4 characters) and nudge its color
# (back to the same value) so the run engine splits it into its own run,
# matching the three-run structure produced by PowerPoint when the text
# was retyped.
$nameSub8 = $tr8.Characters($start8 + 4, 10)
$nameSub8.Font.Color.RGB = 5551484

# ---------------------------------------------------------------------------
# Slides 3, 4, 5, 6: the single-space runs that sit right after "func"/"("/
# "var" keywords had an empty lang="" attribute in the source; PowerPoint
# normalizes this back to lang="en-US" the next time the run is touched.
# ---------------------------------------------------------------------------
function Fix-LangSpace($textRange, $anchorText, $occurrence) {
    $searchStart = 1
    $foundRange = $null
    for ($i = 0; $i -lt $occurrence; $i++) {
        $foundRange = $textRange.Find($anchorText, $searchStart)
        $searchStart = $foundRange.Start + 1
    }
    $spaceStart = $foundRange.Start + $foundRange.Length
    $spaceRun = $textRange.Characters($spaceStart, 1)
    $spaceRun.LanguageID = 1033
    $spaceRun.Text = $spaceRun.Text
}

$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(4).TextFrame.TextRange
Fix-LangSpace $tr3 "func" 1
Fix-LangSpace $tr3 "(" 1
Fix-LangSpace $tr3 "(" 2

$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(4).TextFrame.TextRange
Fix-LangSpace $tr4 "func" 1
Fix-LangSpace $tr4 "func" 2

$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(4).TextFrame.TextRange
Fix-LangSpace $tr5 "var" 1

$tr6b = $s6.Shapes.Item(4).TextFrame.TextRange
Fix-LangSpace $tr6b "var" 1
